$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force the run covering [start, end) to stay/become its own run by
# toggling Bold on then back off across that exact span. Toggling a boolean
# character property on and back off leaves the run's rPr empty (so
# visually/semantically unchanged) but prevents the COM layer from
# re-coalescing it with its neighbours, which is what we need to reproduce
# the multi-run layout the target document wants.
# ---------------------------------------------------------------------------
function Split-Segment($start, $end) {
    $mark = $d.Range($start, $end)
    $mark.Font.Bold = $true
    $mark.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# Locate the paragraph of interest (" and here again in prespecifying
# style:") by its distinctive text rather than hard-coded offsets, so the
# script is resilient to any earlier differences in layout.
# ---------------------------------------------------------------------------
$t = $d.Content.Text
$anchor = " and here again in prespecifying style:"
$base = $t.IndexOf($anchor)

# ---------------------------------------------------------------------------
# 1) Drop the gramStart/gramEnd proofErr pair that wraps "and" by forcing a
#    real text mutation across that span (identical replacement text is a
#    no-op for the engine, so we round-trip through a throwaway marker
#    first) and then re-splitting the resulting merged run into the four
#    pieces the target layout wants: " a" | "nd here again " | "in" | " ".
# ---------------------------------------------------------------------------
$andSpanLen = " and here again in ".Length
$andStart = $base
$andEnd = $base + $andSpanLen

$tmp = $d.Range($andStart, $andEnd)
$tmp.Text = " and here again in #"
$restore = $d.Range($andStart, $andEnd + 1)
$restore.Text = " and here again in "

$p1 = $andStart + 2            # end of " a"
$p2 = $p1 + 14                 # end of "nd here again "
$p3 = $p2 + 2                  # end of "in"
# $andEnd is the end of the trailing " "

Split-Segment $andStart $p1
Split-Segment $p1 $p2
Split-Segment $p2 $p3
Split-Segment $p3 $andEnd

# ---------------------------------------------------------------------------
# 2) "specifying" -> "specified", then split "prespecified" into three runs:
#    "pre" | "specif" | "ied" (spellStart/spellEnd stay put around them).
# ---------------------------------------------------------------------------
$t2 = $d.Content.Text
$specifyingStart = $t2.IndexOf("specifying", $andStart)
$specifyingEnd = $specifyingStart + "specifying".Length
$specifyingRange = $d.Range($specifyingStart, $specifyingEnd)
$specifyingRange.Text = "specified"

$preStart = $specifyingStart - "pre".Length
$specifStart = $specifyingStart
$specifEnd = $specifyingStart + "specif".Length
$iedEnd = $specifyingStart + "specified".Length

Split-Segment $preStart $specifStart
Split-Segment $specifStart $specifEnd
Split-Segment $specifEnd $iedEnd

# ---------------------------------------------------------------------------
# 3) Insert " (red normal)" right after "specified" / before the spellEnd +
#    " style:" text.
# ---------------------------------------------------------------------------
$t3 = $d.Content.Text
$specifiedStart = $t3.IndexOf("specified", $andStart)
$specifiedEnd = $specifiedStart + "specified".Length
$insertPoint = $d.Range($specifiedEnd, $specifiedEnd)
$insertPoint.InsertAfter(" (red normal)")

# ---------------------------------------------------------------------------
# 4) Move the _GoBack bookmark from the trailing empty paragraph to right
#    after "(red normal)" (i.e. right before " style:").
# ---------------------------------------------------------------------------
$t4 = $d.Content.Text
$redNormalEnd = $t4.IndexOf("(red normal)", $andStart) + "(red normal)".Length
$bookmarkRange = $d.Range($redNormalEnd, $redNormalEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
